# Update column G ("K") values on Sheet1, rows 2-68, per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 3
    11 = 0
    12 = 0
    13 = 0
    14 = 2
    15 = 3
    16 = 4
    17 = 3
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    23 = 0
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 0
    35 = 1
    36 = 2
    37 = 0
    38 = 1
    39 = 2
    40 = 1
    41 = 0
    42 = 3
    43 = 1
    44 = 0
    45 = 1
    46 = 0
    47 = 1
    48 = 0
    49 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 2
    54 = 1
    55 = 0
    56 = 0
    57 = 0
    58 = 0
    59 = 1
    60 = 0
    61 = 0
    62 = 0
    63 = 1
    65 = 1
    67 = 1
    68 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
